# Insert a new inventory line (item 22) right after the current item 21 row (row 24).
# The original single line item gets split into two rows: row 24 keeps item #21 but now
# refers to a different product, and the brand-new row 25 picks up the product that used
# to live in row 24, shifting the totals/footer rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new row, then stamp it with row 24's exact formatting (fonts,
# fills, borders, row height, etc.) via PasteSpecial so no new cell styles get minted
# in styles.xml. PasteSpecial does not recreate the source's merged cells, so those are
# re-applied explicitly afterwards.
$ws.Rows("25:25").Insert()
$ws.Range("A24:N24").Copy()
$ws.Range("A25:N25").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("B25:G25").Merge()
$ws.Range("H25:K25").Merge()
$ws.Range("L25:M25").Merge()
$ws.Rows("25:25").RowHeight = 25.5
$ws.Rows("26:26").RowHeight = 25.5

# --- Row 24 (item 21) keeps its row number but now refers to a different product ---
$ws.Range("A24").Value = 21
$ws.Range("B24").Value = "مرطب شفاه لونا جوز هند ابيض"
$ws.Range("H24").Value = "2:0"
$ws.Range("L24").Value = 20
$ws.Range("N24").Value = "1:0"

# --- Row 25 (new item 22) holds the product that used to be item 21 ---
$ws.Range("A25").Value = 22
$ws.Range("B25").Value = "معجون اسنان فلورو بالكولا"
$ws.Range("H25").Value = "3:0"
$ws.Range("L25").Value = 30
$ws.Range("N25").Value = "1:0"

# --- Row 26 (shifted totals row) reflects the re-split quantities ---
$ws.Range("K26").Value = 1533.04
